$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 32.544082
$ws.Range("H2").Value = 97.63224599999999
$ws.Range("I2").Value = 0.621589875979724
$ws.Range("J2").Value = 0.6366365948489335
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.407336
$ws.Range("N2").Value = 0.814672
$ws.Range("Q2").Value = 13.256376185552
$ws.Range("R2").Value = 79.53825711331199
$ws.Range("S2").Value = 0.621589875979724
$ws.Range("T2").Value = 0.6366365948489335

# Row 3
$ws.Range("D3").Value = "MuSCs"
$ws.Range("I3").Value = 0.004665102012661462
$ws.Range("J3").Value = 0.004778029332093849
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.407336
$ws.Range("N3").Value = 0.814672
$ws.Range("Q3").Value = 0.099490595992
$ws.Range("R3").Value = 0.596943575952
$ws.Range("S3").Value = 0.004665102012661462
$ws.Range("T3").Value = 0.004778029332093849

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 6.619872666666667
$ws.Range("H4").Value = 19.859618
$ws.Range("I4").Value = 0.126439142756428
$ws.Range("J4").Value = 0.1294998332673878
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.407336
$ws.Range("N4").Value = 0.814672
$ws.Range("Q4").Value = 2.696512452549333
$ws.Range("R4").Value = 16.179074715296
$ws.Range("S4").Value = 0.126439142756428
$ws.Range("T4").Value = 0.1294998332673878

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 3.71227
$ws.Range("H5").Value = 7.42454
$ws.Range("I5").Value = 0.07090411857072049
$ws.Range("J5").Value = 0.04841365488938666
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.407336
$ws.Range("N5").Value = 0.814672
$ws.Range("Q5").Value = 1.51214121272
$ws.Range("R5").Value = 6.04856485088
$ws.Range("S5").Value = 0.07090411857072049
$ws.Range("T5").Value = 0.04841365488938666

# Row 6
$ws.Range("D6").Value = "MuSCs"
$ws.Range("G6").Value = 9.235725333333333
$ws.Range("H6").Value = 27.707176
$ws.Range("I6").Value = 0.176401760680466
$ws.Range("J6").Value = 0.1806718876621981
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.407336
$ws.Range("N6").Value = 0.814672
$ws.Range("Q6").Value = 3.762043414378666
$ws.Range("R6").Value = 22.572260486272
$ws.Range("S6").Value = 0.176401760680466
$ws.Range("T6").Value = 0.1806718876621981
